$d = $word.ActiveDocument

function Replace-WithLineBreaks {
    param(
        [string]$OldText,
        [string[]]$Parts,
        [bool]$Italic
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text"
    }

    # Clear the matched text but keep the (now empty) run/paragraph in place.
    $rng.Text = ""

    $rPr = ""
    if ($Italic) {
        $rPr = "<w:rPr><w:i/></w:rPr>"
    }

    $runInner = ""
    for ($i = 0; $i -lt $Parts.Length; $i++) {
        if ($i -gt 0) {
            $runInner += "<w:br/>"
        }
        $part = $Parts[$i]
        $escaped = $part -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
        if ($part.StartsWith(" ") -or $part.EndsWith(" ")) {
            $runInner += "<w:t xml:space=`"preserve`">$escaped</w:t>"
        } else {
            $runInner += "<w:t>$escaped</w:t>"
        }
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $rPr + $runInner + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}

# --- Portuguese paragraph ---
$oldPt = "A disciplina consiste no desenvolvimento de um trabalho acadêmico supervisionado por um docente e/ou profissional com título de doutor (orientador), nos seguintes contextos: 1. Desenvolvimento de um projeto de Engenharia, podendo abordar a problemática trabalhada na disciplina LOT2062 Solução de Problemas de Engenharia; 2. Realização de pesquisa científica ou tecnológica inédita, de caráter teórico ou experimental, abordando temas relevantes na área de Engenharia que demandem atualização e síntese de informações."
$partsPt = @(
    "A disciplina consiste no desenvolvimento de um trabalho acadêmico supervisionado por um docente e/ou profissional com título de doutor (orientador), nos seguintes contextos: ",
    "1. Desenvolvimento de um projeto de Engenharia, podendo abordar a problemática trabalhada na disciplina LOT2062 Solução de Problemas de Engenharia; ",
    "2. Realização de pesquisa científica ou tecnológica inédita, de caráter teórico ou experimental, abordando temas relevantes na área de Engenharia que demandem atualização e síntese de informações."
)
Replace-WithLineBreaks $oldPt $partsPt $false

# --- English (italic) paragraph ---
$oldEn = "The course consists of the development of an academic paper supervised by a faculty member and/or professional with a PhD degree (advisor), in the following contexts:1. Development of an Engineering project, which may address the issues covered in the LOT2062 Engineering Problem Solving course;2. Conducting original scientific or technological research, either theoretical or experimental, addressing relevant topics in the field of Engineering that require updating and synthesis of information"
$partsEn = @(
    "The course consists of the development of an academic paper supervised by a faculty member and/or professional with a PhD degree (advisor), in the following contexts:",
    "1. Development of an Engineering project, which may address the issues covered in the LOT2062 Engineering Problem Solving course;",
    "2. Conducting original scientific or technological research, either theoretical or experimental, addressing relevant topics in the field of Engineering that require updating and synthesis of information"
)
Replace-WithLineBreaks $oldEn $partsEn $true
